$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: "T.A. : 2018" -> drop the year, keep just the colon ("'" forces the
# quote-prefix text style so the leading ":" isn't misread as a formula,
# matching the original cell's text formatting)
$ws.Range("J2").Value = "'" + ":"

# Clear the "Guna Membayar" narrative paragraph (merged D11:J12)
$ws.Range("D11").Value = ""

# Clear the "Penerima" signer block (name + NIP)
$ws.Range("A28").Value = ""
$ws.Range("A29").Value = ""

# Clear the "Bendahara Pengeluaran" signer block (name + NIP)
$ws.Range("A39").Value = ""
$ws.Range("A40").Value = ""

# Move the visible selection to the lower signature block
$ws.Range("A42:F43").Select() | Out-Null
